# NN test on dataset-1
# Fill in the S (train MAPE) and T (test MAPE) results for the
# "NN (64,64,64,64,1) (citiesdataset-1.csv)" block on the sheet, plus
# the summary AVERAGE / STDEV.S formulas, and restore the active
# selection/view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @{ Row = 4;  S = 2.3338817159762728; T = 1.844426143127708 },
    @{ Row = 5;  S = 1.8396662312479291; T = 2.898220183133986 },
    @{ Row = 6;  S = 2.6832698447032728; T = 2.5747925879576421 },
    @{ Row = 7;  S = 2.5560328331158351; T = 4.1262811741880876 },
    @{ Row = 8;  S = 2.130629583818636;  T = 2.174855434444773 },
    @{ Row = 9;  S = 2.1367701992429531; T = 2.8717321863135301 },
    @{ Row = 10; S = 1.92050790082374;   T = 6.3663036506270876 },
    @{ Row = 11; S = 2.0474093052567781; T = 2.0402131088053679 },
    @{ Row = 12; S = 1.9040628265116191; T = 2.1729117188903819 },
    @{ Row = 13; S = 2.158731334582503;  T = 3.0445896917428339 },
    @{ Row = 14; S = 2.1653477504282299; T = 2.8508723057229748 },
    @{ Row = 15; S = 1.815639288814815;  T = 2.6667692290704879 },
    @{ Row = 16; S = 2.13641541998236;   T = 1.7513266206440681 },
    @{ Row = 17; S = 2.0677706749864031; T = 3.8699078755347069 },
    @{ Row = 18; S = 1.6608359499970731; T = 8.2983822503209286 },
    @{ Row = 19; S = 1.423836789774076;  T = 7.3366965566606037 },
    @{ Row = 20; S = 2.388303879367482;  T = 1.774771100212797 },
    @{ Row = 21; S = 2.1199140455566901; T = 2.2196718036937231 },
    @{ Row = 22; S = 2.5710529058009661; T = 1.639812421188751 },
    @{ Row = 23; S = 2.1234460015807608; T = 1.4706705227527821 }
)

foreach ($item in $values) {
    $ws.Range("S$($item.Row)").Value = $item.S
    $ws.Range("T$($item.Row)").Value = $item.T
}

$ws.Range("S25").Formula = "=AVERAGE(S4:S23)"
$ws.Range("T25").Formula = "=AVERAGE(T4:T23)"
$ws.Range("S26").Formula = "=STDEV.S(S4:S23)"
$ws.Range("T26").Formula = "=STDEV.S(T4:T23)"

# Restore the view / selection state recorded after the edit.
$ws.Activate()
$ws.Range("AF26").Select()
